# New Test Plans/Try TestCycle
#
# A new automated test cycle was executed against the RMA Complete Flow
# (Issue Credit) scenario. Refresh the RMA numbers / related Salesforce
# record ids captured on the "RMA Details Maintenance Grid" sheet with the
# values produced by the latest run (RMA-1GFW-* -> RMA-R5ZJ-*).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (RMA line 1)
$ws.Range("E2").Value = "RMA-R5ZJ-001"
$ws.Range("F2").Value = "RMA-R5ZJ-1-1"
$ws.Range("J2").Value = "a7s5f000000xKZUAA2"

# Row 3 (RMA line 2)
$ws.Range("E3").Value = "RMA-R5ZJ-002"
$ws.Range("F3").Value = "RMA-R5ZJ-1-2"
$ws.Range("J3").Value = "a7s5f000000xKZVAA2"

# Row 4 (RMA line 3)
$ws.Range("E4").Value = "RMA-R5ZJ-003"
$ws.Range("F4").Value = "RMA-R5ZJ-1-3"
$ws.Range("J4").Value = "a7s5f000000xKZWAA2"

Write-Output "Updated RMA Details Maintenance Grid with new test cycle ids"
